$d = $word.ActiveDocument

function Break-At($findText, $replaceText) {
    $d.Content.Find.Execute($findText, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $replaceText, 2)
}

# --- Paragraph "Programa" (Portuguese) ---
Break-At "Modelos PERT/COM2. Programação" "Modelos PERT/COM^l2. Programação"
Break-At "(branchand-bound).3. Programação Dinâmica3. Métodos" "(branchand-bound).^l3. Programação Dinâmica^l3. Métodos"
Break-At "otimização.4. Modelos e Técnicas" "otimização.^l4. Modelos e Técnicas"

# --- Paragraph "Programa" (English, italic) ---
Break-At "PERT / COM models2. Whole" "PERT / COM models^l2. Whole"
Break-At "Branch-bound algorithm.3. Dynamic Programming3. Heuristic" "Branch-bound algorithm.^l3. Dynamic Programming^l3. Heuristic"
Break-At "optimization problems.4. Forecasting" "optimization problems.^l4. Forecasting"

# --- Paragraph "Bibliografia" ---
Break-At "McGraw-Hill, 2006.2. LACHTERMACHER" "McGraw-Hill, 2006.^l2. LACHTERMACHER"
Break-At "Editora Campus, 2009.3. ANDERSON" "Editora Campus, 2009.^l3. ANDERSON"
Break-At "Publishing, 2000.4. PIZZOLATO" "Publishing, 2000.^l4. PIZZOLATO"
Break-At "LTC Editora, 2009.5. TAHA" "LTC Editora, 2009.^l5. TAHA"
